# Survey2nd add update 4.3 Two_dimensional code
#
# This reproduces three effects visible in the canonical-OOXML diff:
#  1. The "_GoBack" bookmark (Word's automatic "last edit location" marker)
#     moves from its old spot (an empty paragraph after the IoT survey
#     citation) to a new spot in the middle of the CSCWD citation, right
#     after "...cooperative work in desi|gn (CSCWD)...". Word implements
#     this by splitting the run there and inserting the bookmark between
#     the two halves; the old bookmark is removed automatically because
#     "_GoBack" is a singleton bookmark that Word keeps moving around.
#  2. Because of (1), the previously-empty paragraph's bookmarkStart /
#     bookmarkEnd pair disappears (handled automatically by Word when the
#     bookmark is re-added elsewhere).
#  3. A stray straight double-quote immediately after a curly right
#     double-quote (at the very end of the document, after the last
#     citation's URL) is deleted: `...pdf”"` -> `...pdf”`.

$d = $word.ActiveDocument

# --- Change 1 & 2: move the "_GoBack" bookmark -----------------------------
# Locate the split point "...cooperative work in desi|gn (CSCWD)..." via
# Find, collapse the match range to its end, and (re)insert the "_GoBack"
# bookmark there. Word's bookmark model treats "_GoBack" as a single,
# uniquely-named bookmark, so adding it here automatically removes it from
# wherever it previously sat in the document.
$r = $d.Content
$null = $r.Find.Execute(
    "computer supported cooperative work in desi",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)

# --- Change 3: drop the stray straight quote after the closing curly quote -
$r2 = $d.Content
$null = $r2.Find.Execute(
    [char]8221 + '"', $true, $false, $false, $false, $false, $true, 1,
    $false, [char]8221, 2)
